$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 70
$ws.Range("H70").Value = 2574.889
$ws.Range("I70").Value = 846.25
$ws.Range("K70").Value = 2538.75
$ws.Range("M70").Value = -2268.75
# row 73
$ws.Range("H73").Value = 2574.889
$ws.Range("I73").Value = 846.25
$ws.Range("K73").Value = 2538.75
$ws.Range("M73").Value = -1602.75
# row 75
$ws.Range("H75").Value = 28502.334
$ws.Range("I75").Value = 2000
$ws.Range("J75").Value = 33802.8
$ws.Range("K75").Value = 2000
$ws.Range("L75").Value = 33802.8
$ws.Range("M75").Value = -1064
$ws.Range("N75").Value = -35674.8
# row 78
$ws.Range("H78").Value = 28502.334
$ws.Range("I78").Value = 2000
$ws.Range("J78").Value = 33802.8
$ws.Range("K78").Value = 6000
$ws.Range("L78").Value = 101408.4
$ws.Range("M78").Value = -1320
$ws.Range("N78").Value = -110768.4
# row 81
$ws.Range("H81").Value = 41000
$ws.Range("J81").Value = 41000
$ws.Range("L81").Value = 41000
$ws.Range("N81").Value = -42996
# row 84
$ws.Range("H84").Value = 41000
$ws.Range("J84").Value = 41000
$ws.Range("L84").Value = 123000
$ws.Range("N84").Value = -132984
# row 98
$ws.Range("H98").Value = 4068.6316
$ws.Range("I98").Value = 2199.8572
$ws.Range("K98").Value = 2199.8572
$ws.Range("M98").Value = -701.8571999999999
# row 122
$ws.Range("H122").Value = 4068.6316
$ws.Range("I122").Value = 2199.8572
$ws.Range("K122").Value = 6599.571599999999
$ws.Range("M122").Value = -4149.571599999999
# row 124
$ws.Range("H124").Value = 41795
$ws.Range("J124").Value = 41795
$ws.Range("L124").Value = 41795
$ws.Range("N124").Value = -51615
# row 125
$ws.Range("H125").Value = 3961
$ws.Range("I125").Value = 3305.1667
$ws.Range("J125").Value = 4318.727
$ws.Range("K125").Value = 29746.5003
$ws.Range("L125").Value = 38868.543
$ws.Range("M125").Value = -27286.5003
$ws.Range("N125").Value = -43788.543

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 33
$ws.Range("H33").Value = 33005
$ws.Range("J33").Value = 38000
$ws.Range("L33").Value = 38000
$ws.Range("N33").Value = -38658
# row 74
$ws.Range("H74").Value = 5520.409
$ws.Range("I74").Value = 6416.067
$ws.Range("J74").Value = 3601.1428
$ws.Range("K74").Value = 6416.067
$ws.Range("L74").Value = 3601.1428
$ws.Range("M74").Value = -5542.067
$ws.Range("N74").Value = -5349.1428
# row 77
$ws.Range("H77").Value = 5520.409
$ws.Range("I77").Value = 6416.067
$ws.Range("J77").Value = 3601.1428
$ws.Range("K77").Value = 32080.335
$ws.Range("L77").Value = 18005.714
$ws.Range("M77").Value = -27712.335
$ws.Range("N77").Value = -26741.714
# row 122
$ws.Range("H122").Value = 2710
$ws.Range("I122").Value = 1538.7273
$ws.Range("J122").Value = 7004.6665
$ws.Range("K122").Value = 4616.1819
$ws.Range("L122").Value = 21013.9995
$ws.Range("M122").Value = -2166.1819
$ws.Range("N122").Value = -25913.9995
# row 128
$ws.Range("H128").Value = 41824
$ws.Range("J128").Value = 41824
$ws.Range("L128").Value = 41824
$ws.Range("N128").Value = -51784
# row 132
$ws.Range("H132").Value = 2481.5
$ws.Range("I132").Value = 1809.2727
$ws.Range("J132").Value = 4946.3335
$ws.Range("K132").Value = 5427.8181
$ws.Range("L132").Value = 14839.0005
$ws.Range("M132").Value = -2897.8181
$ws.Range("N132").Value = -19899.0005
# row 137
$ws.Range("H137").Value = 39481.54
$ws.Range("J137").Value = 39481.54
$ws.Range("L137").Value = 39481.54
$ws.Range("N137").Value = -49681.54

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 59
$ws.Range("H59").Value = 64474.332
$ws.Range("J59").Value = 64474.332
$ws.Range("L59").Value = 64474.332
$ws.Range("N59").Value = -66168.33199999999
# row 134
$ws.Range("H134").Value = 2833.121
$ws.Range("I134").Value = 1271.6111
$ws.Range("J134").Value = 4706.933
$ws.Range("K134").Value = 3814.8333
$ws.Range("L134").Value = 14120.799
$ws.Range("M134").Value = -1279.8333
$ws.Range("N134").Value = -19190.799
# row 137
$ws.Range("H137").Value = 36586.332
$ws.Range("J137").Value = 38903.6
$ws.Range("L137").Value = 38903.6
$ws.Range("N137").Value = -49103.6

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 5
$ws.Range("H5").Value = 2375.625
$ws.Range("I5").Value = 403
$ws.Range("J5").Value = 3033.1667
$ws.Range("K5").Value = 403
$ws.Range("L5").Value = 3033.1667
$ws.Range("M5").Value = -291
$ws.Range("N5").Value = -3257.1667
# row 31
$ws.Range("H31").Value = 14289117
$ws.Range("I31").Value = 1824.2941
$ws.Range("J31").Value = 27782672
$ws.Range("K31").Value = 1824.2941
$ws.Range("L31").Value = 27782672
$ws.Range("M31").Value = -1529.2941
$ws.Range("N31").Value = -27783262
# row 34
$ws.Range("H34").Value = 14289117
$ws.Range("I34").Value = 1824.2941
$ws.Range("J34").Value = 27782672
$ws.Range("K34").Value = 1824.2941
$ws.Range("L34").Value = 27782672
$ws.Range("M34").Value = -1622.2941
$ws.Range("N34").Value = -27783076
# row 58
$ws.Range("H58").Value = 3614.3635
$ws.Range("I58").Value = 1335.6
$ws.Range("J58").Value = 5513.3335
$ws.Range("K58").Value = 1335.6
$ws.Range("L58").Value = 5513.3335
$ws.Range("M58").Value = -1132.6
$ws.Range("N58").Value = -5919.3335
# row 136
$ws.Range("H136").Value = 3614.3635
$ws.Range("I136").Value = 1335.6
$ws.Range("J136").Value = 5513.3335
$ws.Range("K136").Value = 4006.8
$ws.Range("L136").Value = 16540.0005
$ws.Range("M136").Value = -1456.8
$ws.Range("N136").Value = -21640.0005

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 14
$ws.Range("H14").Value = 746.1667
$ws.Range("I14").Value = 746.1667
$ws.Range("K14").Value = 2238.5001
$ws.Range("M14").Value = -2065.5001
# row 29
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 300
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -623
$ws.Range("N29").ClearContents()
# row 87
$ws.Range("H87").Value = 6925.6
$ws.Range("I87").Value = 3009.3333
$ws.Range("K87").Value = 9027.999899999999
$ws.Range("M87").Value = -7779.999899999999
# row 90
$ws.Range("H90").Value = 6925.6
$ws.Range("I90").Value = 3009.3333
$ws.Range("K90").Value = 27083.9997
$ws.Range("M90").Value = -20843.9997
# row 113
$ws.Range("H113").Value = 720.76086
$ws.Range("I113").Value = 604.1
$ws.Range("J113").Value = 939.5
$ws.Range("K113").Value = 1812.3
$ws.Range("L113").Value = 2818.5
$ws.Range("M113").Value = 357.6999999999998
$ws.Range("N113").Value = -7158.5
# row 131
$ws.Range("H131").Value = 9095366
$ws.Range("I131").Value = 100040220
$ws.Range("J131").Value = 879.5599999999999
$ws.Range("K131").Value = 300120660
$ws.Range("L131").Value = 2638.68
$ws.Range("M131").Value = -300115620
$ws.Range("N131").Value = -12718.68

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 46
$ws.Range("H46").Value = 32303.334
$ws.Range("J46").Value = 32303.334
$ws.Range("L46").Value = 32303.334
$ws.Range("N46").Value = -32615.334
# row 122
$ws.Range("H122").Value = 3681.389
$ws.Range("I122").Value = 1996.091
$ws.Range("K122").Value = 5988.272999999999
$ws.Range("M122").Value = -3538.272999999999
# row 137
$ws.Range("H137").Value = 42655
$ws.Range("J137").Value = 42655
$ws.Range("L137").Value = 42655
$ws.Range("N137").Value = -52855

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 5199.673
$ws.Range("I40").Value = 3455.1428
$ws.Range("J40").Value = 8791.352999999999
$ws.Range("K40").Value = 3455.1428
$ws.Range("L40").Value = 8791.352999999999
$ws.Range("M40").Value = -3319.1428
$ws.Range("N40").Value = -9063.352999999999
# row 46
$ws.Range("H46").Value = 3162.75
$ws.Range("I46").Value = 3733.3333
$ws.Range("K46").Value = 3733.3333
$ws.Range("M46").Value = -3545.3333
# row 55
$ws.Range("H55").Value = 610.1875
$ws.Range("I55").Value = 304.33334
$ws.Range("J55").Value = 1003.4286
$ws.Range("K55").Value = 304.33334
$ws.Range("L55").Value = 1003.4286
$ws.Range("M55").Value = -131.33334
$ws.Range("N55").Value = -1349.4286
# row 58
$ws.Range("H58").Value = 34500
$ws.Range("J58").Value = 34500
$ws.Range("L58").Value = 34500
$ws.Range("N58").Value = -35020
# row 82
$ws.Range("H82").Value = 5156.5186
$ws.Range("I82").Value = 6908.9375
$ws.Range("J82").Value = 2607.5454
$ws.Range("K82").Value = 6908.9375
$ws.Range("L82").Value = 2607.5454
$ws.Range("M82").Value = -6547.9375
$ws.Range("N82").Value = -3329.5454
# row 85
$ws.Range("H85").Value = 5156.5186
$ws.Range("I85").Value = 6908.9375
$ws.Range("J85").Value = 2607.5454
$ws.Range("K85").Value = 6908.9375
$ws.Range("L85").Value = 2607.5454
$ws.Range("M85").Value = -5660.9375
$ws.Range("N85").Value = -5103.5454
# row 122
$ws.Range("H122").Value = 5944.5386
$ws.Range("I122").Value = 3468.5715
$ws.Range("J122").Value = 8833.166999999999
$ws.Range("K122").Value = 10405.7145
$ws.Range("L122").Value = 26499.501
$ws.Range("M122").Value = -7955.7145
$ws.Range("N122").Value = -31399.501
# row 136
$ws.Range("H136").Value = 5858.3335
$ws.Range("I136").Value = 1825
$ws.Range("J136").Value = 7875
$ws.Range("K136").Value = 5475
$ws.Range("L136").Value = 23625
$ws.Range("M136").Value = -2925
$ws.Range("N136").Value = -28725

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 56
$ws.Range("H56").Value = 13016.143
$ws.Range("I56").Value = 1785
$ws.Range("J56").Value = 14888
$ws.Range("K56").Value = 1785
$ws.Range("L56").Value = 14888
$ws.Range("M56").Value = -1071
$ws.Range("N56").Value = -16316
# row 62
$ws.Range("H62").Value = 21851.143
$ws.Range("I62").Value = 8826.333000000001
$ws.Range("K62").Value = 8826.333000000001
$ws.Range("M62").Value = -8202.333000000001
# row 65
$ws.Range("H65").Value = 21851.143
$ws.Range("I65").Value = 8826.333000000001
$ws.Range("K65").Value = 44131.665
$ws.Range("M65").Value = -41011.665
# row 122
$ws.Range("H122").Value = 3016.8
$ws.Range("I122").Value = 1745.5555
$ws.Range("J122").Value = 6285.7144
$ws.Range("K122").Value = 5236.666499999999
$ws.Range("L122").Value = 18857.1432
$ws.Range("M122").Value = -2786.666499999999
$ws.Range("N122").Value = -23757.1432
# row 136
$ws.Range("H136").Value = 2809.158
$ws.Range("I136").Value = 1912.5714
$ws.Range("J136").Value = 3332.1667
$ws.Range("K136").Value = 5737.7142
$ws.Range("L136").Value = 9996.500100000001
$ws.Range("M136").Value = -3187.7142
$ws.Range("N136").Value = -15096.5001
